# Summer 2024 Working Hours - add new log entry (2024-07-11)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 48: Date | Day of week | Hours | (gap) | Notes
$ws.Range("A48").Value = 45484
$ws.Range("A48").NumberFormat = $ws.Range("A47").NumberFormat

$ws.Range("B48").Value = "T"
$ws.Range("B48").NumberFormat = $ws.Range("B47").NumberFormat

$ws.Range("C48").Value = 7

$ws.Range("E48").Value = "working on compiling all code, considering adjusting student grades to either: two different named dfs, or altering the structure of .1 .2 to accommodate other code (would be more work but we can use first attempt data easily)"

# Match row height used for the new wrapped note row
$ws.Rows.Item(48).RowHeight = 90

# Move the active selection to the newly added cell, mirroring the author's last click
$ws.Range("C48").Select() | Out-Null
